$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '327.82'
$r.ClearFormats()

$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = '-1.41%'
$r.ClearFormats()

$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '44.38'
$r.ClearFormats()

$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = '-0.94%'
$r.ClearFormats()

$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = '5.398'
$r.ClearFormats()

$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = '-2.70%'
$r.ClearFormats()

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '0.08364'
$r.ClearFormats()

$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = '0.77%'
$r.ClearFormats()

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '4.426'
$r.ClearFormats()

$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = '-0.17%'
$r.ClearFormats()

$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = '-5.27%'
$r.ClearFormats()

$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '0.9728'
$r.ClearFormats()

$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = '-0.47%'
$r.ClearFormats()

$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = '-4.42%'
$r.ClearFormats()

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '0.1137'
$r.ClearFormats()

$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = '1.50%'
$r.ClearFormats()

$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.1906'
$r.ClearFormats()

$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = '-0.57%'
$r.ClearFormats()

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '0.09675'
$r.ClearFormats()

$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = '-3.97%'
$r.ClearFormats()

$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '0.04605'
$r.ClearFormats()

$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = '-0.65%'
$r.ClearFormats()

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '0.1061'
$r.ClearFormats()

$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = '0.23%'
$r.ClearFormats()

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '0.001291'
$r.ClearFormats()

$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = '2.25%'
$r.ClearFormats()

$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '0.005977'
$r.ClearFormats()

$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = '-0.26%'
$r.ClearFormats()

$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '3.365'
$r.ClearFormats()

$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = '0.14%'
$r.ClearFormats()

$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = '0.10%'
$r.ClearFormats()

$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '8.304'
$r.ClearFormats()

$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = '-19.14%'
$r.ClearFormats()

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '0.1351'
$r.ClearFormats()

$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = '-2.42%'
$r.ClearFormats()

$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = '6.65%'
$r.ClearFormats()

$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '0.04178'
$r.ClearFormats()

$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = '1.64%'
$r.ClearFormats()

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '0.001238'
$r.ClearFormats()

$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = '-4.94%'
$r.ClearFormats()

$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '0.004462'
$r.ClearFormats()

$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = '1.09%'
$r.ClearFormats()

$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '0.0001301'
$r.ClearFormats()

$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = '1.70%'
$r.ClearFormats()

$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '0.0002978'
$r.ClearFormats()

$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = '-20.34%'
$r.ClearFormats()

$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '0.02710'
$r.ClearFormats()

$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = '-3.63%'
$r.ClearFormats()

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '0.05624'
$r.ClearFormats()

$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = '-2.43%'
$r.ClearFormats()

$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '0.007822'
$r.ClearFormats()

$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = '2.34%'
$r.ClearFormats()

$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '0.1414'
$r.ClearFormats()

$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = '-0.97%'
$r.ClearFormats()

$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '0.007289'
$r.ClearFormats()

$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = '-3.45%'
$r.ClearFormats()

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '0.002041'
$r.ClearFormats()

$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = '3.51%'
$r.ClearFormats()

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '0.008722'
$r.ClearFormats()

$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = '8.65%'
$r.ClearFormats()

$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '0.3508'
$r.ClearFormats()

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '0.00006921'
$r.ClearFormats()

$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = '-1.44%'
$r.ClearFormats()

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '0.00000000751'
$r.ClearFormats()

$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = '0.14%'
$r.ClearFormats()

$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '0.003487'
$r.ClearFormats()

$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = '-1.79%'
$r.ClearFormats()

$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '0.003529'
$r.ClearFormats()

$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = '39.81%'
$r.ClearFormats()

$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '0.00002101'
$r.ClearFormats()

$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = '0.14%'
$r.ClearFormats()

$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '0.0002001'
$r.ClearFormats()

$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = '0.14%'
$r.ClearFormats()
